$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 'скважина'
$ws.Cells.Item(2, 3).Value = 0.2282964269090711
$ws.Cells.Item(2, 4).Value = 'добывать скважина, добывать скважина очаг, скважина очаг, нагнетательный скважина'

$ws.Cells.Item(3, 2).Value = 'скв'
$ws.Cells.Item(3, 3).Value = 0.2095315326627219
$ws.Cells.Item(3, 4).Value = 'скважина, добывать скважина, добывать скважина очаг, скважина очаг'

$ws.Cells.Item(4, 2).Value = 'очаг'
$ws.Cells.Item(4, 3).Value = 0.2027955808466008
$ws.Cells.Item(4, 4).Value = 'очаг нагнетание, добывать скважина очаг, скважина очаг, очаг скв'

$ws.Cells.Item(5, 2).Value = 'заводнение'
$ws.Cells.Item(5, 3).Value = 0.1676252261301775
$ws.Cells.Item(5, 4).Value = 'эффективность заводнение, процесс заводнение, оценка эффективность заводнение, заводнение различный'

$ws.Cells.Item(6, 2).Value = 'значение'
$ws.Cells.Item(6, 3).Value = 0.1285799401740974
$ws.Cells.Item(6, 4).Value = 'изменение значение, значение коэффициент, значение очаг, изменение значение коэффициент'

$ws.Cells.Item(7, 2).Value = 'очаг нагнетание'
$ws.Cells.Item(7, 3).Value = 0.128566652808574
$ws.Cells.Item(7, 4).Value = 'очаг нагнетание скв, очаг нагнетание скважина'

$ws.Cells.Item(8, 2).Value = 'изменение значение'
$ws.Cells.Item(8, 3).Value = 0.1257527682755146
$ws.Cells.Item(8, 4).Value = 'изменение значение коэффициент, изменение значение скважина, скважина изменение значение'

$ws.Cells.Item(9, 2).Value = 'коэффициент'
$ws.Cells.Item(9, 3).Value = 0.1071038628702217
$ws.Cells.Item(9, 4).Value = 'значение коэффициент, изменение коэффициент проницаемость, изменение коэффициент, коэффициент проницаемость'

$ws.Cells.Item(10, 2).Value = 'зарифовый мелководье'
$ws.Cells.Item(10, 3).Value = 0.1051908977524696
$ws.Cells.Item(10, 4).Value = 'зона зарифовый мелководье, зарифовый мелководье рифовый'

$ws.Cells.Item(11, 2).Value = 'зарифовый'
$ws.Cells.Item(11, 3).Value = 0.1051908977524696
$ws.Cells.Item(11, 4).Value = 'зарифовый мелководье, зона зарифовый, зона зарифовый мелководье, фациальный зона зарифовый'

$ws.Cells.Item(12, 2).Value = 'эффективность заводнение'
$ws.Cells.Item(12, 3).Value = 0.1004430473651141
$ws.Cells.Item(12, 4).Value = 'оценка эффективность заводнение'

$ws.Cells.Item(13, 2).Value = 'добывать скважина'
$ws.Cells.Item(13, 3).Value = 0.09543863365683752
$ws.Cells.Item(13, 4).Value = 'добывать скважина очаг, коэффициент добывать скважина, фронт добывать скважина'

$ws.Cells.Item(14, 2).Value = 'закачка'
$ws.Cells.Item(14, 3).Value = 0.09518957143501156
$ws.Cells.Item(14, 4).Value = 'накопить закачка, закачка рабочий агент, закачка рабочий, нефть накопить закачка'

$ws.Cells.Item(15, 2).Value = 'добывать'
$ws.Cells.Item(15, 3).Value = 0.09469509279688046
$ws.Cells.Item(15, 4).Value = 'добывать скважина, добывать скважина очаг, скважина добывать, коэффициент добывать'

$ws.Cells.Item(16, 2).Value = 'рифовый'
$ws.Cells.Item(16, 3).Value = 0.08900347769826436
$ws.Cells.Item(16, 4).Value = 'зарифовый мелководье, зарифовый, рифовый гребень, зона зарифовый'

$ws.Cells.Item(17, 2).Value = 'нагнетание'
$ws.Cells.Item(17, 3).Value = 0.08381261306508876
$ws.Cells.Item(17, 4).Value = 'очаг нагнетание, очаг нагнетание скв, нагнетание скв, очаг нагнетание скважина'

$ws.Cells.Item(18, 2).Value = 'значение коэффициент'
$ws.Cells.Item(18, 3).Value = 0.08325549266051276
$ws.Cells.Item(18, 4).Value = 'изменение значение коэффициент, значение коэффициент добывать'

$ws.Cells.Item(19, 2).Value = 'изменение коэффициент проницаемость'
$ws.Cells.Item(19, 3).Value = 0.08181514269636525
$ws.Cells.Item(19, 4).Value = ''

$ws.Cells.Item(20, 2).Value = 'добывать скважина очаг'
$ws.Cells.Item(20, 3).Value = 0.08181514269636525
$ws.Cells.Item(20, 4).Value = ''

$ws.Cells.Item(21, 2).Value = 'рифовый гребень'
$ws.Cells.Item(21, 3).Value = 0.08181514269636525
$ws.Cells.Item(21, 4).Value = ''

$ws.Cells.Item(22, 2).Value = 'мелководье'
$ws.Cells.Item(22, 3).Value = 0.0809907487451607
$ws.Cells.Item(22, 4).Value = 'зарифовый мелководье, зона зарифовый мелководье, мелководье рифовый, зарифовый мелководье рифовый'

$ws.Cells.Item(23, 2).Value = 'фациальный зона'
$ws.Cells.Item(23, 3).Value = 0.0791142023984572
$ws.Cells.Item(23, 4).Value = 'фациальный зона зарифовый'

$ws.Cells.Item(24, 2).Value = 'скважина очаг'
$ws.Cells.Item(24, 3).Value = 0.07812237017286656
$ws.Cells.Item(24, 4).Value = 'добывать скважина очаг, скважина очаг скв'

$ws.Cells.Item(25, 2).Value = 'проницаемость'
$ws.Cells.Item(25, 3).Value = 0.0780524968177891
$ws.Cells.Item(25, 4).Value = 'изменение коэффициент проницаемость, коэффициент проницаемость'

$ws.Cells.Item(26, 2).Value = 'изменение коэффициент'
$ws.Cells.Item(26, 3).Value = 0.07738631893877819
$ws.Cells.Item(26, 4).Value = 'изменение коэффициент проницаемость'

$ws.Cells.Item(27, 2).Value = 'нагнетательный'
$ws.Cells.Item(27, 3).Value = 0.07685286067135809
$ws.Cells.Item(27, 4).Value = 'нагнетательный скважина, нагнетательный скважина добывать, фронт нагнетать нагнетательный'

$ws.Cells.Item(28, 2).Value = 'гребень'
$ws.Cells.Item(28, 3).Value = 0.07525803327461529
$ws.Cells.Item(28, 4).Value = 'рифовый гребень'

$ws.Cells.Item(29, 2).Value = 'озёрный'
$ws.Cells.Item(29, 3).Value = 0.07291769962214874
$ws.Cells.Item(29, 4).Value = 'озёрный месторождение, фм озёрный месторождение, фм озёрный'

$ws.Cells.Item(30, 2).Value = 'фациальный'
$ws.Cells.Item(30, 3).Value = 0.07230820802285964
$ws.Cells.Item(30, 4).Value = 'фациальный зона, фациальный зона зарифовый'

$ws.Cells.Item(31, 2).Value = 'изменение'
$ws.Cells.Item(31, 3).Value = 0.07055046159262758
$ws.Cells.Item(31, 4).Value = 'изменение значение, изменение коэффициент проницаемость, изменение коэффициент, скважина изменение'
